$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the data values for columns I and J (rows 2-12)
$iValues = @(5, 7, 8, 6, 7, 8, 9, 7, 4, 3, 5)
$jValues = @(6, 9, 8, 7, 7, 8, 9, 7, 4, 3, 5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
